$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 12:05"

# 2) Estados Unidos (row 4) - refreshed counters
$ws.Range("B4").Value = 1725808
$ws.Range("C4").Value = 533
$ws.Range("E4").Value = 1145210
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 100625

# 3) Banglades / Bielorrusia swap order (rows 25-26) with refreshed counters
$ws.Range("A25").Value = "Bielorrusia"
$ws.Range("B25").Value = 38956
$ws.Range("C25").Value = 897
$ws.Range("D25").Value = 15923
$ws.Range("E25").Value = 22819
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 214

$ws.Range("A26").Value = "Banglades"
$ws.Range("B26").Value = 38292
$ws.Range("C26").Value = 1541
$ws.Range("D26").Value = 7925
$ws.Range("E26").Value = 29823
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 544

# 4) Colombia / Kuwait swap order (rows 36-37) with refreshed counters
$ws.Range("A36").Value = "Kuwait"
$ws.Range("B36").Value = 23267
$ws.Range("C36").Value = 692
$ws.Range("D36").Value = 7946
$ws.Range("E36").Value = 15146
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 175

$ws.Range("A37").Value = "Colombia"
$ws.Range("B37").Value = 23003
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 5511
$ws.Range("E37").Value = 16716
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 776

# 5) Rumania (row 41) - refreshed counters, no reordering
$ws.Range("B41").Value = 18594
$ws.Range("C41").Value = 165
$ws.Range("D41").Value = 12162
$ws.Range("E41").Value = 5213

# 6) Nigeria / Oman swap order (rows 58-59) with refreshed counters
$ws.Range("A58").Value = "Oman"
$ws.Range("B58").Value = 8373
$ws.Range("C58").Value = 255
$ws.Range("D58").Value = 2177
$ws.Range("E58").Value = 6158
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 38

$ws.Range("A59").Value = "Nigeria"
$ws.Range("B59").Value = 8344
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 2385
$ws.Range("E59").Value = 5710
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 249

# 7) Angola (row 179) - refreshed counters
$ws.Range("B179").Value = 71
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 49
